# Fase 3, informe final y scripts de datos V2
#
# Sprint_Backlog_Actividades_Horas_DFF.xlsx edits:
#   - T-7.2 "Exportar reporte de perdidas a PDF" (row 20) hours raised 3 -> 5
#     (total hours cell C2 recalculates 298 -> 300 automatically)
#   - Selection moved from D34 to D11 (and the view scrolls back up so the
#     scrolled-down top-left anchor is no longer needed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hours estimate for T-7.2
$ws.Range("D20").Value = 5

# Move the active selection to D11 (also resets scroll position)
$ws.Range("D11").Select()
